# Regenerate the "K" column (column G) of save_data with re-computed strike
# counts (s_vals), replacing the old "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K").
$newK = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 2
    11 = 4
    12 = 1
    13 = 1
    14 = 3
    15 = 0
    16 = 1
    17 = 4
    18 = 1
    19 = 3
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 4
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 1
    40 = 1
    41 = 1
    42 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
